# Added test data for invalid login / Created InvalidLogin
#
# Starting point: a single sheet "test1" with A1 = "akshara".
# Target: two sheets - "ValidLogin" (renamed from the original sheet) and a
# new "InvalidLogin" sheet, each holding a small UserName/Password login
# table, with InvalidLogin left as the active/selected sheet.

$wb = $excel.ActiveWorkbook

# Rename the existing (only) sheet to "ValidLogin".
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "ValidLogin"

# Duplicate it right after itself to get "InvalidLogin" - this keeps the
# same worksheet formatting/namespaces as the original sheet instead of
# picking up a brand-new blank-sheet template.
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "InvalidLogin"

# Valid-login test data.
$ws1.Range("A1").Value = "UserName"
$ws1.Range("B1").Value = "Password"
$ws1.Range("A2").Value = "admin"
$ws1.Range("B2").Value = "manager"

# Invalid-login test data.
$ws2.Range("A1").Value = "UserName"
$ws2.Range("B1").Value = "Password"
$ws2.Range("A2").Value = "abcd"
$ws2.Range("B2").Value = "xyz"

# Match the saved selections/active sheet from the source workbook.
$ws1.Range("B1").Select()
$ws2.Range("E6").Select()
